# Port skip and take to GetFlights
# Adds a new "GetFlightsRows" worksheet (after the existing "GetFlights" sheet)
# containing the StartRow/EndRow/ExpectedCount/ExpectedRows test matrix used
# to exercise the new skip/take paging logic.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "GetFlightsRows"

# Everything on this sheet is stored as text (same convention as the other
# test-data sheets in this workbook).
$ws.Range("A1:G10").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "Flights"
$ws.Range("C1").Value = "StartRow"
$ws.Range("D1").Value = "EndRow"
$ws.Range("F1").Value = "ExpectedCount"
$ws.Range("G1").Value = "ExpectedRows"

# Row 2
$ws.Range("A2").Value = "3"
$ws.Range("C2").Value = "-1"
$ws.Range("D2").Value = "-1"
$ws.Range("F2").Value = "3"
$ws.Range("G2").Value = "1,2,3"

# Row 3
$ws.Range("A3").Value = "3"
$ws.Range("C3").Value = "0"
$ws.Range("D3").Value = "0"
$ws.Range("F3").Value = "1"
$ws.Range("G3").Value = "1"

# Row 4
$ws.Range("A4").Value = "3"
$ws.Range("C4").Value = "1"
$ws.Range("D4").Value = "1"
$ws.Range("F4").Value = "1"
$ws.Range("G4").Value = "2"

# Row 5
$ws.Range("A5").Value = "3"
$ws.Range("C5").Value = "2"
$ws.Range("D5").Value = "2"
$ws.Range("F5").Value = "1"
$ws.Range("G5").Value = "3"

# Row 6
$ws.Range("A6").Value = "3"
$ws.Range("C6").Value = "0"
$ws.Range("D6").Value = "1"
$ws.Range("F6").Value = "2"
$ws.Range("G6").Value = "1,2"

# Row 7
$ws.Range("A7").Value = "3"
$ws.Range("C7").Value = "1"
$ws.Range("D7").Value = "2"
$ws.Range("F7").Value = "2"
$ws.Range("G7").Value = "2,3"

# Row 8
$ws.Range("A8").Value = "3"
$ws.Range("C8").Value = "-1"
$ws.Range("D8").Value = "1"
$ws.Range("F8").Value = "2"
$ws.Range("G8").Value = "1,2"

# Row 9
$ws.Range("A9").Value = "3"
$ws.Range("C9").Value = "1"
$ws.Range("D9").Value = "-1"
$ws.Range("F9").Value = "2"
$ws.Range("G9").Value = "2,3"

# Row 10
$ws.Range("A10").Value = "3"
$ws.Range("C10").Value = "0"
$ws.Range("D10").Value = "999"
$ws.Range("F10").Value = "3"
$ws.Range("G10").Value = "1,2,3"

# Selection / view bookkeeping so the new sheet matches the rest of the
# workbook's sheetView conventions and becomes the active tab.
$ws.Range("A1:G10").Select() | Out-Null
